$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("x – 1;", $true, $false, $false, $false, $false, $true, 1, $false, "x - 1;", 2)
